$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 20 (their content is merged into rows 2-5)
$ws.Range("A6:A20").EntireRow.Delete()

# Update the remaining data rows with the combined tuple-style strings
$ws.Range("A2").Value = '(''Dragon Fodder'', [''{1}{R}'', ''Sorcery'', ''Create two 1/1 red Goblin creature tokens.''])'
$ws.Range("A3").Value = '("Dragonlord''s Servant", [''{1}{R}'', ''Creature — Goblin Shaman'', ''Dragon spells you cast cost {1} less to cast.'', ''1/3''])'
$ws.Range("A4").Value = '(''Evolving Wilds'', [''Land'', ''{T}, Sacrifice Evolving Wilds: Search your library for a basic land card, put it onto the battlefield tapped, then shuffle your library.''])'
$ws.Range("A5").Value = '(''Foe-Razer Regent'', [''{5}{G}{G}'', ''Creature — Dragon'', ''Flying'', ''When Foe-Razer Regent enters the battlefield, you may have it fight target creature you don’t control.'', ''Whenever a creature you control fights, put two +1/+1 counters on it at the beginning of the next end step.'', ''4/5''])'

